$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "98.359.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.83%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.487.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.64%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "665.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.15%  "

$ws.Range("E7").Value = "  +5.92%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.427"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.27%  "

$ws.Range("E9").Value = "  +2.94%  "

$ws.Range("E10").Value = "  +0.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.487.32"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.50%  "

$ws.Range("E13").Value = "  +0.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "98.081.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.78%  "

$ws.Range("E16").Value = "  +1.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.141.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.04%  "

$ws.Range("E18").Value = "  +3.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.489.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +10.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.527"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.99%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "522.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.53%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000203"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.68%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "98.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.675.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.00%  "

$ws.Range("E30").Value = "  +11.64%  "

$ws.Range("E31").Value = "  +13.48%  "

$ws.Range("E32").Value = "  -1.77%  "

$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("E34").Value = "  -0.86%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.596"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.07%  "

$ws.Range("E37").Value = "  +0.19%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.87%  "

$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.02%  "

$ws.Range("E40").Value = "  +4.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "525.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.50%  "

$ws.Range("E42").Value = "  +0.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.913"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.16%  "

$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "24.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.91%  "

$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0434"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.67%  "

$ws.Range("E47").Value = "  +2.88%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.02%  "

$ws.Range("E49").Value = "  -1.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.74%  "
